$d = $word.ActiveDocument

$pairs = @(
    @("367×9=", "442×2="),
    @("422×3=", "146×2="),
    @("223×3=", "397×3="),
    @("800×4=", "867×6="),
    @("571×5=", "303×3="),
    @("183×9=", "860×5="),
    @("568×4=", "839×9="),
    @("666×6=", "573×6="),
    @("293×9=", "653×8="),
    @("944×4=", "369×8="),
    @("281×5=", "819×8="),
    @("988×2=", "274×3="),
    @("623×7=", "329×3="),
    @("903×4=", "822×9="),
    @("294×4=", "104×7="),
    @("927×4=", "374×8="),
    @("523×4=", "345×9="),
    @("840×4=", "127×3="),
    @("310×4=", "750×2="),
    @("415×9=", "637×2="),
    @("288×2=", "609×5="),
    @("467×3=", "991×6="),
    @("268×4=", "700×2="),
    @("265×3=", "203×5="),
    @("545×9=", "874×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
